$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "23.851.16"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "1.639.13"
$ws.Range("E3").Value = "  -1.50%  "

$ws.Range("D4").Value = "1.0000"
$ws.Range("E4").Value = "  -0.03%  "

$ws.Range("D5").Value = "308.66"
$ws.Range("E5").Value = "  -0.44%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").Value = "0.3858"
$ws.Range("E7").Value = "  -1.23%  "

$ws.Range("D8").Value = "0.3810"
$ws.Range("E8").Value = "  -1.68%  "

$ws.Range("D9").Value = "50.70"
$ws.Range("E9").Value = "  -1.88%  "

$ws.Range("D10").Value = "1.322"
$ws.Range("E10").Value = "  -3.73%  "

$ws.Range("D11").Value = "0.9999"
$ws.Range("E11").Value = "  -0.03%  "

$ws.Range("D12").Value = "0.08362"
$ws.Range("E12").Value = "  -1.74%  "

$ws.Range("D13").Value = "23.67"
$ws.Range("E13").Value = "  -2.03%  "

$ws.Range("D14").Value = "6.942"
$ws.Range("E14").Value = "  -4.09%  "

$ws.Range("D15").Value = "7.759"
$ws.Range("E15").Value = "  -3.49%  "

$ws.Range("E16").Value = "  -1.41%  "

$ws.Range("D17").Value = "1.622.79"
$ws.Range("E17").Value = "  -2.37%  "

$ws.Range("D18").Value = "93.47"
$ws.Range("E18").Value = "  -1.48%  "

$ws.Range("D19").Value = "0.06929"
$ws.Range("E19").Value = "  -1.00%  "

$ws.Range("D20").Value = "19.36"
$ws.Range("E20").Value = "  -3.27%  "

$ws.Range("D21").Value = "6.842"
$ws.Range("E21").Value = "  -2.32%  "

$ws.Range("D22").Value = "1.000"
$ws.Range("E22").Value = "  -0.03%  "

$ws.Range("D23").Value = "13.48"
$ws.Range("E23").Value = "  -1.86%  "

$ws.Range("D24").Value = "23.854.96"
$ws.Range("E24").Value = "  -0.79%  "

$ws.Range("D25").Value = "2.430"
$ws.Range("E25").Value = "  -3.19%  "

$ws.Range("D26").Value = "2.869"
$ws.Range("E26").Value = "  -9.01%  "

$ws.Range("D27").Value = "21.77"
$ws.Range("E27").Value = "  -2.55%  "

$ws.Range("D28").Value = "152.62"
$ws.Range("E28").Value = "  -0.78%  "

$ws.Range("D29").Value = "5.459"
$ws.Range("E29").Value = "  +2.58%  "

$ws.Range("D30").Value = "136.42"
$ws.Range("E30").Value = "  -3.76%  "

$ws.Range("D31").Value = "7.779"
$ws.Range("E31").Value = "  -0.80%  "

$ws.Range("D32").Value = "2.478"
$ws.Range("E32").Value = "  -0.16%  "

$ws.Range("D33").Value = "1.837.61"
$ws.Range("E33").Value = "  -0.31%  "

$ws.Range("D34").Value = "0.07920"

$ws.Range("D35").Value = "0.9773"
$ws.Range("E35").Value = "  -7.70%  "

$ws.Range("D36").Value = "0.02878"
$ws.Range("E36").Value = "  -4.89%  "

$ws.Range("D37").Value = "6.554"
$ws.Range("E37").Value = "  -2.69%  "

$ws.Range("D38").Value = "0.2646"
$ws.Range("E38").Value = "  -2.71%  "

$ws.Range("D39").Value = "10.43"
$ws.Range("E39").Value = "  -7.78%  "

$ws.Range("E40").Value = "  -1.20%  "

$ws.Range("D41").Value = "0.7462"
$ws.Range("E41").Value = "  -2.12%  "

$ws.Range("B42").Value = "Aptos"
$ws.Range("C42").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D42").Value = "13.21"
$ws.Range("E42").Value = "  -3.71%  "

$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").Value = "1.411"
$ws.Range("E43").Value = "  -1.44%  "

$ws.Range("D44").Value = "16.47"
$ws.Range("E44").Value = "  -0.74%  "

$ws.Range("D45").Value = "0.6859"
$ws.Range("E45").Value = "  -2.67%  "

$ws.Range("D46").Value = "2.401"
$ws.Range("E46").Value = "  -4.55%  "

$ws.Range("D47").Value = "4.064"
$ws.Range("E47").Value = "  -0.96%  "

$ws.Range("D48").Value = "0.9997"
$ws.Range("E48").Value = "  -0.06%  "

$ws.Range("D49").Value = "0.08215"
$ws.Range("E49").Value = "  -1.41%  "

$ws.Range("D50").Value = "133.76"
$ws.Range("E50").Value = "  -1.39%  "

$ws.Range("D51").Value = "1.216"
$ws.Range("E51").Value = "  -2.28%  "
